# "feat: true dataset compression test"
# Adds a new "Sheet2" (after the existing "Sheet1") holding a compression
# benchmark table (step / my / protobuf / json byte sizes), makes it the
# active sheet, and leaves Sheet1 as before (just no longer the selected tab).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after Sheet1 and rename it.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row.
$ws2.Cells.Item(1, 1).Value = "step"
$ws2.Cells.Item(1, 2).Value = "my"
$ws2.Cells.Item(1, 3).Value = "protobuf"
$ws2.Cells.Item(1, 4).Value = "json"

# Data rows (step, my-size, protobuf-size, json-size).
$data = @(
    @(5, 78618394, 143080950, 324611758),
    @(10, 80710521, 143080950, 324611758),
    @(15, 82491112, 143080950, 324611758),
    @(20, 83839352, 143080950, 324611758),
    @(25, 85375420, 143080950, 324611758),
    @(30, 86617185, 143080950, 324611758),
    @(35, 88006772, 143080950, 324611758),
    @(40, 89105098, 143080950, 324611758),
    @(45, 90438945, 143080950, 324611758),
    @(50, 91517267, 143080950, 324611758),
    @(100, 99054933, 143080950, 324611758),
    @(150, 102067017, 143080950, 324611758),
    @(200, 103920808, 143080950, 324611758),
    @(250, 105124270, 143080950, 324611758),
    @(300, 106015527, 143080950, 324611758),
    @(350, 106757016, 143080950, 324611758),
    @(400, 107342854, 143080950, 324611758),
    @(450, 107808411, 143080950, 324611758),
    @(500, 108211898, 143080950, 324611758),
    @(550, 108539643, 143080950, 324611758),
    @(600, 108854963, 143080950, 324611758),
    @(650, 109255003, 143080950, 324611758),
    @(700, 109378730, 143080950, 324611758),
    @(750, 109831735, 143080950, 324611758),
    @(800, 109892320, 143080950, 324611758),
    @(850, 110158661, 143080950, 324611758),
    @(900, 110191041, 143080950, 324611758),
    @(950, 110203215, 143080950, 324611758),
    @(1000, 110397692, 143080950, 324611758)
)

$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt 4; $c++) {
        $ws2.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# Row 31 is intentionally left blank; row 32 carries the trailing note.
$ws2.Cells.Item(32, 1).Value = "Process finished with exit code 0"

# Match the recorded selection / active sheet state.
[void]$ws2.Range("H12").Select()
[void]$ws2.Activate()
